$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("M14").Value = "INTRPT"
$ws.Range("M15").Value = "RETINT"
$ws.Range("M16").Value = "NOINT"
$ws.Range("M17").Value = "ENINT"

$ws.Range("M18").Select()
